# Insert a new record row at row 587, shifting existing rows 587:660 down to 588:661.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 587 (pushes old 587..660 -> 588..661)
$ws.Rows.Item(587).Insert()

# Populate the newly inserted row 587 with the new record's values.
$ws.Cells.Item(587, 1).Value = 10
$ws.Cells.Item(587, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(587, 3).Value = "La Araucanía"
$ws.Cells.Item(587, 4).Value = 44918
$ws.Cells.Item(587, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(587, 5).Value = 9
$ws.Cells.Item(587, 6).Value = 100112043
$ws.Cells.Item(587, 7).Value = "Pepino ensalada"
$ws.Cells.Item(587, 8).Value = "Sin especificar"
$ws.Cells.Item(587, 9).Value = "Primera"
$ws.Cells.Item(587, 10).Value = 125
$ws.Cells.Item(587, 11).Value = 15000
$ws.Cells.Item(587, 12).Value = 15000
$ws.Cells.Item(587, 13).Value = 15000
$ws.Cells.Item(587, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(587, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(587, 16).Value = 250
$ws.Cells.Item(587, 17).Value = 60
$ws.Cells.Item(587, 18).Value = "Hortaliza"
